# Script update and archiving.
# Adds a "Second batch mkii" tracking column (S) and a "Count" column (T)
# that records whether each split-GAL4 line's data was Deleted / Loaded /
# Dropped as empty / Merged, plus the resulting row-count, and a running
# total in T15. Also updates a couple of stale "ToDo" notes in column K.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells (row 1) ---
$ws.Range("S1").Value = "Second batch mkii"
$ws.Range("T1").Value = "Count"

# --- Row 2 ---
$ws.Range("R2").Value = "Deleted"
$ws.Range("S2").Value = "Loaded"
$ws.Range("T2").Value = 188

# --- Row 3 ---
$ws.Range("R3").Value = "Deleted"
$ws.Range("S3").Value = "Loaded"
$ws.Range("T3").Value = 40

# --- Row 4 ---
$ws.Range("R4").Value = "Deleted"
$ws.Range("S4").Value = "Dropped as empty"
$ws.Range("T4").Value = 0

# --- Row 5 ---
$ws.Range("R5").Value = "Deleted"
$ws.Range("S5").Value = "Loaded"
$ws.Range("T5").Value = 8

# --- Row 6 ---
$ws.Range("K6").Value = "Merge with pub ver"
$ws.Range("R6").Value = "Merge"
$ws.Range("S6").Value = "Merged with Dolan2019"
$ws.Range("T6").Value = 0

# --- Row 7 ---
$ws.Range("R7").Value = "Deleted"
$ws.Range("S7").Value = "Loaded"
$ws.Range("T7").Value = 840

# --- Row 8 ---
$ws.Range("R8").Value = "Deleted"
$ws.Range("S8").Value = "Loaded"
$ws.Range("T8").Value = 4

# --- Row 9 ---
$ws.Range("R9").Value = "Deleted"
$ws.Range("S9").Value = "Loaded"
$ws.Range("T9").Value = 8

# --- Row 10 ---
$ws.Range("R10").Value = "Deleted"
$ws.Range("S10").Value = "Loaded"
$ws.Range("T10").Value = 126

# --- Row 11 ---
$ws.Range("K11").Value = "Merge with pre ver"
$ws.Range("R11").Value = "Deleted"
$ws.Range("S11").Value = "Loaded"
$ws.Range("T11").Value = 329

# --- Row 12 ---
$ws.Range("R12").Value = "Deleted"
$ws.Range("S12").Value = "Loaded"
$ws.Range("T12").Value = 256

# --- Row 13: a standalone note in column S ---
$ws.Range("S13").Value = "20x and 63x only when left_dorsal, ventral, right_dorsal tiles present"

# --- Row 15: running total of the new Count column ---
$ws.Range("T15").Formula = "=SUM(T2:T12)"

# Keep the active selection in line with the newly extended used range.
$ws.Range("T16").Select()
